$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.761.04"
$ws.Range("D2").NumberFormat = "General"
$ws.Range("E2").Value = "  -1.11%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.033.72"
$ws.Range("D3").NumberFormat = "General"
$ws.Range("E3").Value = "  -1.47%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "227.40"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  -1.30%  "
$ws.Range("E6").Value = "  -0.71%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.02"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = "  +3.19%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  -0.20%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0816"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  +0.69%  "
$ws.Range("E11").Value = "  +0.09%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.61"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = "  -0.23%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.332.32"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = "  -1.57%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.04"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  +1.23%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.759"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = "  +0.65%  "
$ws.Range("E16").Value = "  -1.44%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.025.77"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = "  -1.87%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.748.93"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = "  -0.91%  "
$ws.Range("E19").Value = "  -2.05%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.81"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = "  -0.14%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0823"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  -1.41%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "225.57"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = "  +0.26%  "
$ws.Range("E23").Value = "  +0.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.36"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  -3.44%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.21"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  -2.28%  "
$ws.Range("E26").Value = "  -0.49%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "165.11"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  -0.40%  "
$ws.Range("E28").Value = "  -4.22%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.90"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  -1.02%  "
$ws.Range("E30").Value = "  -7.30%  "
$ws.Range("E31").Value = "  +1.26%  "
$ws.Range("E32").Value = "  -2.84%  "
$ws.Range("E33").Value = "  +4.14%  "
$ws.Range("E35").Value = "  -3.12%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.41"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = "  +6.28%  "
$ws.Range("E38").Value = "  -2.01%  "
$ws.Range("E39").Value = "  -0.07%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.541.90"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = "  +4.09%  "
$ws.Range("B41").Value = "InjectiveProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "16.97"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  +0.54%  "
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0216"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = "  -1.14%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "96.83"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = "  -1.76%  "
$ws.Range("E44").Value = "  -1.62%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0920"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  -2.70%  "
$ws.Range("E46").Value = "  -1.66%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.91"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  -4.81%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.00"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.96"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = "  -0.51%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.13"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  +0.26%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.223.11"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  -1.52%  "
